# CancelBillingTests.xlsx
# Replace the old "execute query / check value" DB-verification steps with
# the new portal keyword-driven steps (login + service-center "add new"
# flow), and make the header row stand out in bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1 - navigate to portal (args no longer needed)
$ws.Range("C2").Value = "navigate to portal"
$ws.Range("D2").Value = $null
$ws.Range("I2").Value = $null

# Step 2 - log in to portal, with user/pass args
$ws.Range("C3").Value = "log in to portal"
$ws.Range("D3").Value = "QAGENERIC"
$ws.Range("E3").Value = "QA!generic1"

# Step 3 - navigate to service center update section
$ws.Range("C4").Value = "navigate to service center update secion"
$ws.Range("D4").Value = $null
$ws.Range("I4").Value = $null

# Step 4 - click add new
$ws.Range("C5").Value = "service update - click add new"
$ws.Range("D5").Value = $null
$ws.Range("E5").Value = $null

# Step 5 - scroll to custpro area
$ws.Range("C6").Value = "service update add new - scroll to custpro aread"

# Step 6 - select cancel rebill checkbox
$ws.Range("C7").Value = "service update add new - select cancel rebill checkbox"

# Step 7 (new row) - click import button
$ws.Range("A8").Value = "Y"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "service update add new - click import button"

# Bold the header row
$ws.Range("A1:K1").Font.Bold = $true

# Keyword/Argument 1 columns hold different text now, so resize them
$ws.Columns("C:D").AutoFit()

# Leave the selection where the author left it
$ws.Range("C7").Select()
